# ListaVecinos.xlsx edit:
#  - Insert a new row between the current row 63 ("79") and row 64 ("82 / 83"),
#    splitting the combined "82 / 83" entry into two separate rows:
#      new row 64: control 82, Estela Retta, 1163620357, 115 nº 1364
#      row 65 (was row 64): control 83, "-", "-", "-"
#  - Existing rows 65/66 shift down to 66/67
#  - A new trailing row 68 is appended containing a single ellipsis ("…") in A68
#  - Final selection becomes D66

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the trailing ellipsis row first, so it becomes the first newly-added shared
# string (matching the original authoring order). This row is beyond the shifted
# block below, so it does not need to move afterwards.
$ws.Range("A68").Value2 = [char]0x2026

# Shift rows 64-66 down to 65-67 using Range-to-Range copy (preserves existing cell
# styles exactly and does not fabricate new style/cellXf entries the way
# Rows.Insert() does on this engine).
$ws.Range("A66:D66").Copy($ws.Range("A67:D67"))
$ws.Range("A65:D65").Copy($ws.Range("A66:D66"))
$ws.Range("A64:D64").Copy($ws.Range("A65:D65"))

# Row 65 (previously row 64) used to hold the combined "82 / 83" label in column A;
# now it only represents control number 83, the rest of the row ("-") is unchanged.
$ws.Range("A65").Value2 = 83

# Populate the new row 64 with the split-out "82" entry
$ws.Range("A64").Value2 = 82
$ws.Range("B64").Value2 = "Estela Retta"
$ws.Range("C64").Value2 = 1163620357
$ws.Range("D64").Value2 = "115 n" + [char]0x00BA + " 1364"

# Leave the final selection on D66, matching where editing ended
$ws.Range("D66").Select() | Out-Null
